$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# This run immediately follows a hyperlink's closing boundary. Replacing the
# full match (starting exactly at the run boundary) via Find/Replace or by
# assigning .Text on a range that begins right at that boundary causes the
# new text to incorrectly inherit the adjoining hyperlink run's formatting.
# Work around it by leaving the first character (the leading space, which is
# unchanged anyway) untouched and only replacing from the second character
# onward, which keeps the replaced text inside the original (correctly
# formatted) run.
function Replace-TextAfterHyperlinkBoundary($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Start = $rng.Start + 1
    $rng.Text = $replace.Substring(1)
}

Replace-Text "英语" "英語"
Replace-TextAfterHyperlinkBoundary " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語"
Replace-Text "简介" "簡介"
Replace-Text "一封发送给目标国家中未通过我们验证流程的合作伙伴的电子邮件。 将通过 customer.io 发送" "發送給目標國家中那些文件未通過我們驗證流程的合作夥伴的電子郵件。 將通過 customer.io 發送"
Replace-Text "目标受众" "目標受眾"
Replace-Text "提交了错误/不完整文件的邀请合作伙伴" "提交錯誤/不完整文檔的被邀請合作夥伴"
Replace-Text "主题行" "主題行"
Replace-Text "[事件名称]" "[事件名稱]"
Replace-Text " — 文档验证失败 " " — 文件驗證失敗 "
Replace-Text "啊哦！ 文件无法验证" "啊哦！ 文檔無法驗證"
Replace-Text "[合作伙伴姓名]" "[合作夥伴姓名]"
Replace-Text "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "很遺憾地通知您，您的文檔未通過驗證流程，因為我們發現以下問題： "
Replace-Text "疫苗接种证书副本" "您的疫苗接種證明副本"
Replace-Text ": 文件不清楚" ": 文檔不清楚"
Replace-Text "[文件 2]" "[文檔 2]"
Replace-Text ": [problem]" ": [問題]"
Replace-Text "请在 " "請在 "
Replace-Text " 之前重新提交上述文件，以便我们进行必要的安排。" " 之前重新提交上述文檔，以便我們進行必要的安排。"
Replace-Text "如有任何疑问，请通过 " "如有任何疑問，請通過 "
Replace-Text "[电子邮件地址]" "[電子郵件地址]"
Replace-Text "[WHATSAPP 号码]" "[WHATSAPP 號碼]"
Replace-Text " (WhatsApp) 联系您的区域经理 " " (WhatsApp) 聯繫您的區域經理, "
Replace-Text "[NAME]" "[姓名]"
Replace-Text " 。 " "。 "
